# Add daily power records: fill in the still-blank last table row (row 20)
# and append a brand-new row (row 21) to the "comforter_cda_table" Excel table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- Row 20 already existed inside the table (with calculated D/E/F formulas
#     but blank Date/Start Time/End Time). Fill in the missing inputs. ---
$ws.Cells.Item(20, 1).Value = 43347              # Date       -> 2018-09-04
$ws.Cells.Item(20, 2).Value = 0.74236111111111114 # Start Time -> 17:49:00
$ws.Cells.Item(20, 3).Value = 0.99930555555555556 # End Time   -> 23:59:00

# --- Grow the table by one row so it now spans A1:F21. ---
$newRow = $lo.ListRows.Add()

# --- New row 21: Date + Start Time filled in; End Time left blank. ---
$ws.Cells.Item(21, 1).Value = 43348   # Date       -> 2018-09-05
$ws.Cells.Item(21, 2).Value = 0       # Start Time -> 00:00:00

# The engine doesn't auto-propagate the table's calculated-column formulas
# to the newly added row, so set them explicitly (matching the other rows).
$ws.Cells.Item(21, 4).Formula = "=(C21-B21)* 1440"
$ws.Cells.Item(21, 5).Formula = "=IF(C21>B21, (C21-B21)*1440, (B21-C21)*1440)"
$ws.Cells.Item(21, 6).Formula = "=ABS((C21-B21)*1440)"

# Column A ("Date") needs to be a touch wider to fit the new, longer date
# string ("Wednesday, September 05, 2018").
$ws.Columns.Item(1).ColumnWidth = 27.85

# Reflect where the user's selection ended up after entering the new data.
[void]$ws.Range("C21").Select()
